$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B6 value from "Autre" to "Chauffage"
$ws.Range("B6").Value = "Chauffage"

# Move the active selection from B9 to B10 (cursor position side-effect)
$ws.Range("B10").Select()
